$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in column H, row 1, matching the
# formatting (bold, bordered, centered) already used by the other
# header cells (e.g. G1).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# Add the corresponding data value for row 2.
$ws.Range("H2").Value = 1
